# Fgf22-Fgfrl1.xlsx -- refresh LR-pair scores with new TPM-based values
# and add the "ECs" sending-cluster rows (new data rows 2-4; existing
# "MuSCs" sending-cluster rows shift down to rows 5-7 with updated values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Fgf22/Fgfrl1)
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Fgf22"
$ws.Cells.Item(2,3).Value = "Fgfrl1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.09830299999999999
$ws.Cells.Item(2,8).Value = 0.294909
$ws.Cells.Item(2,9).Value = 0.3722197960868258
$ws.Cells.Item(2,10).Value = 0.3722197960868259
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.6537256666666666
$ws.Cells.Item(2,14).Value = 1.961177
$ws.Cells.Item(2,15).Value = 0.08287237534104652
$ws.Cells.Item(2,16).Value = 0.08287237534104651
$ws.Cells.Item(2,17).Value = 0.06426319421033333
$ws.Cells.Item(2,18).Value = 0.578368747893
$ws.Cells.Item(2,19).Value = 0.03084673865067523
$ws.Cells.Item(2,20).Value = 0.03084673865067523

# Row 3: ECs -> FAPs (Fgf22/Fgfrl1)
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Fgf22"
$ws.Cells.Item(3,3).Value = "Fgfrl1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.09830299999999999
$ws.Cells.Item(3,8).Value = 0.294909
$ws.Cells.Item(3,9).Value = 0.3722197960868258
$ws.Cells.Item(3,10).Value = 0.3722197960868259
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 5.625751333333334
$ws.Cells.Item(3,14).Value = 16.877254
$ws.Cells.Item(3,15).Value = 0.7131728182689164
$ws.Cells.Item(3,16).Value = 0.7131728182689163
$ws.Cells.Item(3,17).Value = 0.5530282333206666
$ws.Cells.Item(3,18).Value = 4.977254099886
$ws.Cells.Item(3,19).Value = 0.265457040990723
$ws.Cells.Item(3,20).Value = 0.265457040990723

# Row 4: ECs -> MuSCs (Fgf22/Fgfrl1)
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Fgf22"
$ws.Cells.Item(4,3).Value = "Fgfrl1"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.09830299999999999
$ws.Cells.Item(4,8).Value = 0.294909
$ws.Cells.Item(4,9).Value = 0.3722197960868258
$ws.Cells.Item(4,10).Value = 0.3722197960868259
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.608865333333333
$ws.Cells.Item(4,14).Value = 4.826596
$ws.Cells.Item(4,15).Value = 0.2039548063900371
$ws.Cells.Item(4,16).Value = 0.2039548063900371
$ws.Cells.Item(4,17).Value = 0.1581562888626666
$ws.Cells.Item(4,18).Value = 1.423406599764
$ws.Cells.Item(4,19).Value = 0.07591601644542766
$ws.Cells.Item(4,20).Value = 0.07591601644542766

# Row 5: MuSCs -> ECs (Fgf22/Fgfrl1)
$ws.Cells.Item(5,1).Value = "MuSCs"
$ws.Cells.Item(5,2).Value = "Fgf22"
$ws.Cells.Item(5,3).Value = "Fgfrl1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.1657963333333334
$ws.Cells.Item(5,8).Value = 0.497389
$ws.Cells.Item(5,9).Value = 0.6277802039131741
$ws.Cells.Item(5,10).Value = 0.6277802039131741
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.6537256666666666
$ws.Cells.Item(5,14).Value = 1.961177
$ws.Cells.Item(5,15).Value = 0.08287237534104652
$ws.Cells.Item(5,16).Value = 0.08287237534104651
$ws.Cells.Item(5,17).Value = 0.1083853185392222
$ws.Cells.Item(5,18).Value = 0.975467866853
$ws.Cells.Item(5,19).Value = 0.05202563669037129
$ws.Cells.Item(5,20).Value = 0.05202563669037128

# Row 6: MuSCs -> FAPs (Fgf22/Fgfrl1)
$ws.Cells.Item(6,1).Value = "MuSCs"
$ws.Cells.Item(6,2).Value = "Fgf22"
$ws.Cells.Item(6,3).Value = "Fgfrl1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 0.6666666666666666
$ws.Cells.Item(6,7).Value = 0.1657963333333334
$ws.Cells.Item(6,8).Value = 0.497389
$ws.Cells.Item(6,9).Value = 0.6277802039131741
$ws.Cells.Item(6,10).Value = 0.6277802039131741
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 5.625751333333334
$ws.Cells.Item(6,14).Value = 16.877254
$ws.Cells.Item(6,15).Value = 0.7131728182689164
$ws.Cells.Item(6,16).Value = 0.7131728182689163
$ws.Cells.Item(6,17).Value = 0.9327289433117779
$ws.Cells.Item(6,18).Value = 8.394560489806
$ws.Cells.Item(6,19).Value = 0.4477157772781934
$ws.Cells.Item(6,20).Value = 0.4477157772781933

# Row 7: MuSCs -> MuSCs (Fgf22/Fgfrl1)
$ws.Cells.Item(7,1).Value = "MuSCs"
$ws.Cells.Item(7,2).Value = "Fgf22"
$ws.Cells.Item(7,3).Value = "Fgfrl1"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 0.6666666666666666
$ws.Cells.Item(7,7).Value = 0.1657963333333334
$ws.Cells.Item(7,8).Value = 0.497389
$ws.Cells.Item(7,9).Value = 0.6277802039131741
$ws.Cells.Item(7,10).Value = 0.6277802039131741
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.608865333333333
$ws.Cells.Item(7,14).Value = 4.826596
$ws.Cells.Item(7,15).Value = 0.2039548063900371
$ws.Cells.Item(7,16).Value = 0.2039548063900371
$ws.Cells.Item(7,17).Value = 0.2667439730937778
$ws.Cells.Item(7,18).Value = 2.400695757844
$ws.Cells.Item(7,19).Value = 0.1280387899446094
$ws.Cells.Item(7,20).Value = 0.1280387899446094
